$wb = $excel.ActiveWorkbook

# --- Sheet "TestCases" ---
$ws1 = $wb.Worksheets.Item("TestCases")

# Row 3: Test2/N -> PaymentOptionsTest/Y
$ws1.Range("A3").Value = "PaymentOptionsTest"
$ws1.Range("B3").Value = "Y"

# Column A got wider to fit the new, longer test name
$ws1.Columns.Item(1).ColumnWidth = 22.5

# Selection left where the edit happened
$ws1.Range("E15").Select()

# --- Sheet "Data" ---
$ws2 = $wb.Worksheets.Item("Data")

# Test1 block: flip a couple of Runmode flags
$ws2.Range("A3").Value = "N"
$ws2.Range("A5").Value = "N"

# Rename the "Test2" (stock purchase) section into "PaymentOptionsTest"
# and give the header a plain (non-bold) red fill instead of the bold one.
$ws2.Range("A7").Value = "PaymentOptionsTest"
$ws2.Range("A7").Font.Bold = $false
$ws2.Range("A7").Interior.Color = 255

# New, simplified column headers for the section
$ws2.Range("D8").Value = "Col1"
$ws2.Range("E8").Value = "Col2"
$ws2.Range("F8").Value = "Col3"
$ws2.Range("G8").Value = "Col4"

# Replace the old 3-row stock-purchase data with a single simple data row
$ws2.Range("A9").Value = "Y"
$ws2.Range("C9").Value = "Mozilla"
$ws2.Range("D9:G9").ClearFormats()
$ws2.Range("D9").Value = "C11"
$ws2.Range("E9").Value = "C11"
$ws2.Range("F9").Value = "C35"
$ws2.Range("G9").Value = "C35"

# The old block had 3 data rows, the new one only has 1 -> drop the extra two,
# shifting ViewBillTest/Test4/Test5 blocks up by two rows.
$ws2.Rows("10:11").Delete()

# Column A widened slightly
$ws2.Columns.Item(1).ColumnWidth = 19.67

# Selection left where the edit happened
$ws2.Range("K19").Select()
